$d = $word.ActiveDocument

# 1. Replace the placeholder "Matheus... " run with the full name + RGM.
$find = $d.Content.Find
$find.Execute("Matheus... ", $true, $false, $false, $false, $false, $true, 1, $false, `
    "Matheus Gomes Santos " + [char]0x2013 + " RGM: 26136805", 2)

# 2. Move the lone "_GoBack" bookmark so it sits right after the run we
#    just inserted (instead of its old spot in the empty paragraph near
#    "Link das bibliotecas").
$r = $d.Content
$r.Find.Execute("RGM: 26136805", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)

$old = $d.Bookmarks("_GoBack")
$old.Delete()
$d.Bookmarks.Add("_GoBack", $r)
